$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B8"), "https://uibank.uipath.com/welcome") | Out-Null
$ws.Range("B8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B9"), "https://uibank.uipath.com/accounts") | Out-Null
$ws.Range("B9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B13"), "mailto:sophiasoaconnolly@gmail.com") | Out-Null
$ws.Range("B13").Style = "Hyperlink"
